# "Generate Report for Handback"
# Update status text for the second file (2bb65dd6-...) from "Ready for
# handoff" to "Handback transform failed" across the Overview, zh-cn and
# de-de sheets, and record the handback-transform error detail message for
# each locale, widening the "Error Detail" column to fit the new text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Overview sheet: zh-cn / de-de status columns for row 3
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Per-locale sheets: Status column (C) for row 3
$wsZhCn.Range("C3").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# Per-locale sheets: Error Detail column (P) for row 3
$wsZhCn.Range("P3").Value = "Handback file name: nazvj2k4.rmc is different with handoff file name: 2bb65dd6-be15-4d3e-809a-85941ac2b4f8.f79f984b0d8174995604c58c377eae0297506f83.zh-cn."
$wsDeDe.Range("P3").Value = "Handback file name: nazvj2k4.rmc is different with handoff file name: 2bb65dd6-be15-4d3e-809a-85941ac2b4f8.f79f984b0d8174995604c58c377eae0297506f83.de-de."

# Widen column P (Error Detail) on both locale sheets so the longer message
# fits (stored column width ends up at 40 characters).
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
